$d = $word.ActiveDocument

# The template placeholders {TESTO_1} .. {TESTO_8} (all caps) should become
# the lower-cased {testo_1} .. {testo_8}. Use Find/Replace with MatchCase
# so we only touch the exact upper-case tokens.
for ($i = 1; $i -le 8; $i++) {
    $old = "{TESTO_$i}"
    $new = "{testo_$i}"
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}
